# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-44, replacing the previous Strike#-derived values.
$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 3
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 1
    26 = 2
    27 = 1
    28 = 2
    29 = 2
    30 = 0
    31 = 2
    32 = 1
    33 = 3
    34 = 0
    35 = 1
    36 = 1
    37 = 0
    38 = 3
    39 = 1
    40 = 1
    41 = 2
    42 = 3
    43 = 2
    44 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
